# Append the new rows (6-18) of Name/Dependant/Wives/Projects data to Sheet1.
# Column layout (per existing header row): A=Name, B=Dependant, C=Wives, D=Projects

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Hafsa",  2, 7, 4),
    @("Alice",  5, 4, 11),
    @("John",   3, 7, 9),
    @("Jasmin", 8, 8, 9),
    @("eve",    4, 4, 4),
    @("smith",  3, 2, 3),
    @("james",  5, 5, 5),
    @("winny",  6, 2, 3),
    @("khan",   5, 4, 4),
    @("shams",  5, 7, 7),
    @("aliyah", 7, 1, 1),
    @("danny",  5, 3, 3),
    @("suzy",   5, 5, 3)
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

$ws.Range("A18").Select()
